$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells so numeric-looking values
# (e.g. "1.00", "594.89") are stored as text, matching the source data,
# rather than being auto-converted to numbers by Excel.
$dCells = @("D2","D3","D4","D5","D6","D7","D9","D10","D11","D12","D13","D14","D15","D16","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values (prices, volumes, and re-ordered coin rows)
$ws.Range("D2").Value = '65.054.76'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").Value = '3.546.43'
$ws.Range("E3").Value = '  +3.79%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '594.89'
$ws.Range("E5").Value = '  +2.49%  '
$ws.Range("D6").Value = '137.35'
$ws.Range("E6").Value = '  +2.36%  '
$ws.Range("D7").Value = '3.546.63'
$ws.Range("E7").Value = '  +3.71%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").Value = '0.493'
$ws.Range("E9").Value = '  +2.61%  '
$ws.Range("D10").Value = '0.123'
$ws.Range("E10").Value = '  +3.25%  '
$ws.Range("D11").Value = '6.96'
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").Value = '0.385'
$ws.Range("E12").Value = '  +3.37%  '
$ws.Range("D13").Value = '4.149.88'
$ws.Range("E13").Value = '  +3.87%  '
$ws.Range("D14").Value = '0.0000182'
$ws.Range("E14").Value = '  +3.31%  '
$ws.Range("D15").Value = '27.22'
$ws.Range("E15").Value = '  +4.91%  '
$ws.Range("D16").Value = '3.549.90'
$ws.Range("E16").Value = '  +3.64%  '
$ws.Range("E17").Value = '  +1.45%  '
$ws.Range("D18").Value = '64.900.16'
$ws.Range("E18").Value = '  +0.57%  '
$ws.Range("D19").Value = '10.15'
$ws.Range("E19").Value = '  +7.77%  '
$ws.Range("D20").Value = '5.83'
$ws.Range("E20").Value = '  +1.95%  '
$ws.Range("D21").Value = '14.28'
$ws.Range("E21").Value = '  +6.30%  '
$ws.Range("D22").Value = '389.79'
$ws.Range("E22").Value = '  +3.09%  '
$ws.Range("D23").Value = '0.575'
$ws.Range("E23").Value = '  +6.69%  '
$ws.Range("D24").Value = '3.684.63'
$ws.Range("E24").Value = '  +3.69%  '
$ws.Range("D25").Value = '73.88'
$ws.Range("E25").Value = '  +3.68%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").Value = '0.0000114'
$ws.Range("E27").Value = '  +10.26%  '
$ws.Range("D28").Value = '7.74'
$ws.Range("E28").Value = '  +8.59%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E30").Value = '  +5.15%  '
$ws.Range("D31").Value = '8.23'
$ws.Range("E31").Value = '  +3.55%  '
$ws.Range("D32").Value = '3.556.80'
$ws.Range("E32").Value = '  +3.67%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = '1.38'
$ws.Range("E33").Value = '  +18.11%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = '23.85'
$ws.Range("E35").Value = '  +4.16%  '
$ws.Range("D36").Value = '0.144'
$ws.Range("E36").Value = '  +3.31%  '
$ws.Range("D37").Value = '169.14'
$ws.Range("E37").Value = '  +0.61%  '
$ws.Range("D38").Value = '1.56'
$ws.Range("E38").Value = '  +9.16%  '
$ws.Range("D39").Value = '6.90'
$ws.Range("E39").Value = '  +4.19%  '
$ws.Range("D40").Value = '4.96'
$ws.Range("E40").Value = '  +9.59%  '
$ws.Range("D41").Value = '0.0804'
$ws.Range("E41").Value = '  +7.27%  '
$ws.Range("D42").Value = '0.822'
$ws.Range("E42").Value = '  +1.99%  '
$ws.Range("D43").Value = '26.59'
$ws.Range("E43").Value = '  +18.60%  '
$ws.Range("D44").Value = '42.78'
$ws.Range("E44").Value = '  +2.58%  '
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").Value = '1.22'
$ws.Range("E46").Value = '  +10.10%  '
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").Value = '4.43'
$ws.Range("E47").Value = '  +4.49%  '
$ws.Range("D48").Value = '1.67'
$ws.Range("E48").Value = '  +4.94%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.443.75'
$ws.Range("E49").Value = '  +12.84%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '6.88'
$ws.Range("E50").Value = '  +7.34%  '
$ws.Range("D51").Value = '302.84'
$ws.Range("E51").Value = '  +12.63%  '

# Remove the temporary text-number-format styling so the cells keep
# their original (unstyled) appearance while remaining text values.
foreach ($addr in $dCells) {
    $ws.Range($addr).ClearFormats()
}
